# This workbook tracks daily "Rabanito" (radish) price records for the
# "Vega Central Mapocho de Santiago" market. Two new weekly price records
# need to be inserted into the existing data block (rows 222-329), each as
# a brand new row, pushing the subsequent rows down. The new rows keep the
# same constant fields (Mercado ID, Mercado, Region, Codreg, Categoria ID,
# Categoria, Variedad, Calidad, Unidad de comercializacion, Kg o Unidades,
# Clasificacion) as every other row in this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RabanitoRow($RowNum, $Fecha, $Volumen, $PrecioMinimo, $PrecioMaximo, $PrecioPromedio, $Origen, $PrecioKg) {
    $ws.Cells.Item($RowNum, 1).Value = 9
    $ws.Cells.Item($RowNum, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($RowNum, 3).Value = "Metropolitana"
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 5).Value = 13
    $ws.Cells.Item($RowNum, 6).Value = 300000001
    $ws.Cells.Item($RowNum, 7).Value = "Rabanito"
    $ws.Cells.Item($RowNum, 8).Value = "Sin especificar"
    $ws.Cells.Item($RowNum, 9).Value = "Primera"
    $ws.Cells.Item($RowNum, 10).Value = $Volumen
    $ws.Cells.Item($RowNum, 11).Value = $PrecioMinimo
    $ws.Cells.Item($RowNum, 12).Value = $PrecioMaximo
    $ws.Cells.Item($RowNum, 13).Value = $PrecioPromedio
    $ws.Cells.Item($RowNum, 14).Value = "`$/cien unidades (volumen en unidades)"
    $ws.Cells.Item($RowNum, 15).Value = $Origen
    $ws.Cells.Item($RowNum, 16).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 17).Value = 100
    $ws.Cells.Item($RowNum, 18).Value = "Hortaliza"
}

# Insert the first new record before the existing row 222, pushing the
# whole data block (rows 222-329) down by one row.
$ws.Rows(222).Insert()
Set-RabanitoRow 222 44846 5000 4000 4000 4000 "Provincia de Chacabuco" 40

# Insert the second new record before what is now row 231 (the former
# row 230), pushing the remaining tail of the data block down by one more
# row.
$ws.Rows(231).Insert()
Set-RabanitoRow 231 44845 5000 4000 4000 4000 "Provincia de Chacabuco" 40
